$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data block (rows 33-34), pushing the
# existing rows (old 33 onward) down by two.
$ws.Rows("33:34").Insert()

# New row 33: Asterix / 1a (cosecha lavada) reading dated 2022-02-09
$ws.Range("A33").Value = 1
$ws.Range("B33").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C33").Value = "Arica y Parinacota"
$ws.Range("D33").Value = 44601
$ws.Range("E33").Value = 15
$ws.Range("F33").Value = 100114001
$ws.Range("G33").Value = "Papa"
$ws.Range("H33").Value = "Asterix"
$ws.Range("I33").Value = "1a (cosecha lavada)"
$ws.Range("J33").Value = 1000
$ws.Range("K33").Value = 11000
$ws.Range("L33").Value = 12000
$ws.Range("M33").Value = 11500
$ws.Range("N33").Value = "`$/saco 25 kilos"
$ws.Range("O33").Value = "Región de Los Lagos"
$ws.Range("P33").Value = 460
$ws.Range("Q33").Value = 25
$ws.Range("R33").Value = "Hortaliza"

# New row 34: Patagonia / 1a (cosecha) reading, also dated 2022-02-09
$ws.Range("A34").Value = 1
$ws.Range("B34").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C34").Value = "Arica y Parinacota"
$ws.Range("D34").Value = 44601
$ws.Range("E34").Value = 15
$ws.Range("F34").Value = 100114001
$ws.Range("G34").Value = "Papa"
$ws.Range("H34").Value = "Patagonia"
$ws.Range("I34").Value = "1a (cosecha)"
$ws.Range("J34").Value = 1000
$ws.Range("K34").Value = 9000
$ws.Range("L34").Value = 10000
$ws.Range("M34").Value = 9500
$ws.Range("N34").Value = "`$/saco 25 kilos"
$ws.Range("O34").Value = "Provincia de Melipilla"
$ws.Range("P34").Value = 380
$ws.Range("Q34").Value = 25
$ws.Range("R34").Value = "Hortaliza"

# The old "Cardinal / Región de Coquimbo" reading (originally row 38, now
# pushed down to row 40 by the insert above) is removed entirely.
$ws.Rows("40:40").Delete()
